$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet has headers in row 1 and data rows 2..94.
# Columns: A=code, B=name, C=status, D=codeforiati:group-code, E=codeforiati:group-name
# This edit swaps the contents of columns D and E (including the header row)
# so that D becomes "group-name" and E becomes "group-code".

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
